$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" - update the handoff generation timestamps for the
# 2eaeadba-8695-4808-b297-150735bf2676.md entry across the Overview, zh-cn and de-de sheets.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the 2eaeadba row (row 5)
$wsOverview.Range("G5").Value = "2016-08-30 18:52:35"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the 2eaeadba row (row 5)
$wsZhCn.Range("H5").Value = "2016-08-30 18:52:30"

# de-de sheet: "Latest Handoff Datetime" column (H) for the 2eaeadba row (row 5)
$wsDeDe.Range("H5").Value = "2016-08-30 18:52:35"
